# \n processing added in the annotation strings!
#
# Annotations sheet: the BottomLeft / BottomRight annotation rows swap
# their text+color content. BottomLeft now carries the Legend text
# (rewritten to use literal "\n" separators instead of " / "), and
# BottomRight now carries the "Plot / Map data" credit text (which used
# to live in BottomLeft). Blank row 7 gets its formatting cleared.
#
# Instructions sheet: a new row documents that literal "\n" can be used
# inside annotation strings to force a line break.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Annotations sheet
# ---------------------------------------------------------------------
$annotations = $wb.Worksheets.Item("Annotations")

# Row 5 ("BottomLeft") used to show the map-credit text; it now shows
# the legend text (with literal \n line breaks), no longer wrapped.
$annotations.Cells.Item(5, 2).Value = "Legend:\nGreen=Long Waves\nBlue=Medium Waves"
$annotations.Cells.Item(5, 2).WrapText = $false
$annotations.Cells.Item(5, 4).Value = "black"

# Row 6 ("BottomRight") used to show the legend text; it now shows the
# map-credit text, which is wrapped.
$annotations.Cells.Item(6, 2).Value = "Plot: OsmMarker by Ynovo\nMap data: OpenStreetMap"
$annotations.Cells.Item(6, 2).WrapText = $true
$annotations.Cells.Item(6, 4).Value = "blue"

# Row 7 is blank; clear its leftover wrap formatting and make sure C7/D7
# carry the same (empty / default) formatting as B7.
$annotations.Cells.Item(7, 2).WrapText = $false
$annotations.Cells.Item(7, 3).Value = ""
$annotations.Cells.Item(7, 4).Value = ""

# ---------------------------------------------------------------------
# Instructions sheet
# ---------------------------------------------------------------------
$instructions = $wb.Worksheets.Item("Instructions")
$instructions.Cells.Item(10, 1).Value = "Possibility to use \n in strings to change line"

# ---------------------------------------------------------------------
# Selections: keep "Annotations" the active/selected sheet (as before),
# but move its active cell to B4; Instructions' active cell moves to A10.
# ---------------------------------------------------------------------
$instructions.Activate()
$instructions.Range("A10").Select()

$annotations.Activate()
$annotations.Range("B4").Select()
